# Apply the authored edit to AFC.xlsx:
#  - Rename sheet "TDC" to "Tableau Disjonctif Complet"
#  - Update the active selection on that sheet from G11 to E16

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("TDC")
$ws.Name = "Tableau Disjonctif Complet"

$ws.Activate()
[void]$ws.Range("E16").Select()
